$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values with multiple
# "." separators (e.g. "41.754.88") or plain decimals (e.g. "231.14") are
# not auto-coerced into numbers by Excel's type inference on assignment.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '41.754.88'
$ws.Range('E2').Value = '  +5.45%  '
$ws.Range('D3').Value = '2.225.10'
$ws.Range('E3').Value = '  +2.89%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '231.14'
$ws.Range('E5').Value = '  +1.92%  '
$ws.Range('D6').Value = '0.623'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D7').Value = '60.63'
$ws.Range('E7').Value = '  -3.19%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +2.95%  '
$ws.Range('D10').Value = '58.83'
$ws.Range('E10').Value = '  +0.74%  '
$ws.Range('E11').Value = '  +4.74%  '
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('D13').Value = '2.553.36'
$ws.Range('E13').Value = '  +2.83%  '
$ws.Range('D14').Value = '15.64'
$ws.Range('E14').Value = '  -1.30%  '
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').Value = '0.799'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('E17').Value = '  +1.64%  '
$ws.Range('D18').Value = '2.226.58'
$ws.Range('E18').Value = '  +2.96%  '
$ws.Range('D19').Value = '41.687.95'
$ws.Range('E19').Value = '  +5.39%  '
$ws.Range('D20').Value = '72.73'
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('D23').Value = '250.08'
$ws.Range('E23').Value = '  +9.80%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').Value = '9.62'
$ws.Range('E27').Value = '  +2.09%  '
$ws.Range('E28').Value = '  +3.12%  '
$ws.Range('D29').Value = '167.50'
$ws.Range('E29').Value = '  -1.65%  '
$ws.Range('D30').Value = '19.91'
$ws.Range('E30').Value = '  +1.44%  '
$ws.Range('E31').Value = '  -1.39%  '
$ws.Range('E32').Value = '  -1.93%  '
$ws.Range('D33').Value = '0.122'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').Value = '4.94'
$ws.Range('E34').Value = '  +5.24%  '
$ws.Range('E35').Value = '  +2.29%  '
$ws.Range('D36').Value = '0.0624'
$ws.Range('E36').Value = '  +1.16%  '
$ws.Range('D37').Value = '6.62'
$ws.Range('E37').Value = '  -4.84%  '
$ws.Range('E38').Value = '  -3.52%  '
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('D40').Value = '0.000254'
$ws.Range('E40').Value = '  +28.23%  '
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('E42').Value = '  +5.65%  '
$ws.Range('E43').Value = '  -2.40%  '
$ws.Range('D44').Value = '8.63'
$ws.Range('D45').Value = '0.0980'
$ws.Range('E45').Value = '  +7.17%  '
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').Value = '98.55'
$ws.Range('E47').Value = '  -3.66%  '
$ws.Range('D48').Value = '1.467.73'
$ws.Range('E48').Value = '  -2.95%  '
$ws.Range('D49').Value = '16.49'
$ws.Range('E49').Value = '  -6.63%  '
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').Value = '51.93'
$ws.Range('E51').Value = '  +7.24%  '

# Restore the default (General) style on the price column so no stray
# number-format / style metadata is left behind on these cells.
$priceRange.Style = "Normal"
